$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 73: Green / Electric vehicles / Electric vehicles
$ws.Range("A73").Value = "Green"
$ws.Range("B73").Value = "Electric vehicles"
$ws.Range("C73").Value = "Electric vehicles"

# New row 74: Brown / long description (no Description/column C value)
$ws.Range("A74").Value = "Brown"
$ws.Range("B74").Value = "companies engaged in exploration & production, refining & marketing, and storage & transportation of oil & gas and coal & consumable fuels. It also includes companies that offer oil & gas equipment and services."

# Highlight the newly added rows with a yellow fill, matching the
# actual used cells in each row (row 74 has no C column entry).
$ws.Range("A73:C73").Interior.Color = 65535
$ws.Range("A74:B74").Interior.Color = 65535

# Match the final selection state: whole row 74 selected, active cell A74.
$ws.Range("A74").EntireRow.Select() | Out-Null
